$wb = $excel.ActiveWorkbook

# Sheets that contain the affected rows: "展览" and "全部类型"
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1324
    $ws.Range("F3").Value = 1778
    $ws.Range("F4").Value = 109
    $ws.Range("F6").Value = 6278
    $ws.Range("F7").Value = 128
}
